# plantilla-base-flujo.xlsx: "Detail base routines added, missing fields"
#
# 1) Rename Hoja1 -> Base Flujo, Hoja2 -> Base Detalle
# 2) Populate "Base Detalle" (sheet2) with the same header row layout as
#    "Base Flujo" (sheet1) plus a bunch of new detail columns, and make it
#    the active/selected sheet (Base Flujo was active before).

$wb = $excel.ActiveWorkbook

$wsFlujo   = $wb.Worksheets.Item(1)
$wsDetalle = $wb.Worksheets.Item(2)

$wsFlujo.Name   = "Base Flujo"
$wsDetalle.Name = "Base Detalle"

# --- Base Detalle header row -------------------------------------------------

$headers = @(
    "Punto de Venta",        # A
    "Dirección",              # B
    "Comuna",                 # C
    "Zona",                   # D
    "Fecha",                  # E
    "Encuestador",             # F
    "Persona Nº",             # G
    "Marca",                  # H
    "Modelo",                 # I
    "Tipo de Contrato",       # J
    "Modalidad de Equipo",    # K
    "Otras: Tarjeta",         # L
    "Otras: Chip",            # M
    "Otras: Accesorios",      # N
    "Calificación`nPlan",     # O
    "SEMANA",                 # P
    "MES",                    # Q
    "TIPO PDV",               # R
    "Operador",                # S
    "Calificación Equipo",    # T
    "Recarga Express",        # U
    "Carac Adicional",        # V
    "Total Ventas",           # W
    "Cambios Portabilidad",   # X
    "Razones Portabilidad"    # Y
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsDetalle.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Base formatting for the whole header row: reuse the "Normal 2" style
# (bold white Arial 10 on the red fill, wrap text) already used by Base
# Flujo's header row, cell by cell so the engine's stylesheet writer
# reuses the existing font/fill records instead of minting new ones.
for ($col = 1; $col -le $headers.Length; $col++) {
    $wsFlujo.Range("A1").Copy()
    $wsDetalle.Cells.Item(1, $col).PasteSpecial(-4122)
}

# E1 (Fecha) keeps the date-number-format variant, like Base Flujo's E1.
$wsFlujo.Range("E1").Copy()
$wsDetalle.Range("E1").PasteSpecial(-4122)

# Left-align the descriptive columns: A,B,D,F..O
$leftAddrs = @("A1","B1","D1","F1","G1","H1","I1","J1","K1","L1","M1","N1","O1")
foreach ($addr in $leftAddrs) {
    $wsDetalle.Range($addr).HorizontalAlignment = -4131
}

# Centered "totals" column (still on the original red fill).
$wsDetalle.Range("W1").VerticalAlignment = -4108
$wsDetalle.Range("W1").HorizontalAlignment = -4108

# Centered "portability" columns, with a slightly darker fill variant.
foreach ($addr in @("X1","Y1")) {
    $wsDetalle.Range($addr).VerticalAlignment = -4108
    $wsDetalle.Range($addr).HorizontalAlignment = -4108
    $wsDetalle.Range($addr).Interior.ColorIndex = 3
}

$wsDetalle.Rows.Item(1).RowHeight = 39

$wsDetalle.Range("B19").Select()
$wsDetalle.Activate()
